$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -2.578184724194803
$ws.Range("D2").Value = 0.0144367625963806

$ws.Range("C3").Value = -1.260504182493561
$ws.Range("D3").Value = 0.2160698029907064

$ws.Range("C4").Value = -0.775234558306254
$ws.Range("D4").Value = 0.443559775714621

$ws.Range("C5").Value = -0.08191624050633897
$ws.Range("D5").Value = 0.9351936650984594

$ws.Range("C6").Value = 1.243339391685659
$ws.Range("D6").Value = 0.2222482252978311

$ws.Range("C7").Value = 1.649449676781618
$ws.Range("D7").Value = 0.1082659088089868

$ws.Range("C8").Value = 2.766396732371846
$ws.Range("D8").Value = 0.009098627491888056

$ws.Range("C9").Value = 0.6775427581358235
$ws.Range("D9").Value = 0.5026464416077001

$ws.Range("C10").Value = 0.8627006967195744
$ws.Range("D10").Value = 0.3943466422042712

$ws.Range("C11").Value = 0.4877470573379163
$ws.Range("D11").Value = 0.6288595266659827
